$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 7000
$ws.Range("I116").Value = 6000
$ws.Range("K116").Value = 6000
$ws.Range("M116").Value = -2558

# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 22567.285
$ws.Range("I132").Value = 3500.05
$ws.Range("K132").Value = 10500.15
$ws.Range("M132").Value = -7970.150000000001

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 3061.4546
$ws.Range("I137").Value = 2768.1
$ws.Range("K137").Value = 8304.299999999999
$ws.Range("M137").Value = -5754.299999999999

$ws = $wb.Worksheets.Item("ARM")
# Row 11 (Leve Item ID 3767)
$ws.Range("H11").Value = 5002501.5
$ws.Range("I11").Value = 5002501.5
$ws.Range("K11").Value = 5002501.5
$ws.Range("M11").Value = -5002357.5

# Row 19 (Leve Item ID 3550)
$ws.Range("H19").Value = 3014.2
$ws.Range("I19").Value = 1004
$ws.Range("J19").Value = 4354.3335
$ws.Range("K19").Value = 1004
$ws.Range("L19").Value = 4354.3335
$ws.Range("M19").Value = -775
$ws.Range("N19").Value = -4812.3335

$ws = $wb.Worksheets.Item("BSM")
# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 36985.832
$ws.Range("I99").Value = 36985.832
$ws.Range("K99").Value = 36985.832
$ws.Range("M99").Value = -35487.832

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 1767.6444
$ws.Range("I31").Value = 1678.4706
$ws.Range("J31").Value = 2043.2727
$ws.Range("K31").Value = 1678.4706
$ws.Range("L31").Value = 2043.2727
$ws.Range("M31").Value = -1383.4706
$ws.Range("N31").Value = -2633.2727

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 1767.6444
$ws.Range("I34").Value = 1678.4706
$ws.Range("J34").Value = 2043.2727
$ws.Range("K34").Value = 1678.4706
$ws.Range("L34").Value = 2043.2727
$ws.Range("M34").Value = -1476.4706
$ws.Range("N34").Value = -2447.2727

# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 4396.8
$ws.Range("J58").Value = 4665
$ws.Range("L58").Value = 4665
$ws.Range("N58").Value = -5071

# Row 105 (Leve Item ID 19928)
$ws.Range("H105").Value = 2573.25
$ws.Range("I105").Value = 2717.0625
$ws.Range("K105").Value = 2717.0625
$ws.Range("M105").Value = -970.0625

# Row 122 (Leve Item ID 36196)
$ws.Range("H122").Value = 447825.22
$ws.Range("I122").Value = 929634.6
$ws.Range("K122").Value = 2788903.8
$ws.Range("M122").Value = -2786453.8

# Row 132 (Leve Item ID 44019)
$ws.Range("N132").ClearContents()
$ws.Range("H132").Value = 3346.75
$ws.Range("I132").Value = 3346.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10040.25
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7510.25

# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 4396.8
$ws.Range("J136").Value = 4665
$ws.Range("L136").Value = 13995
$ws.Range("N136").Value = -19095

$ws = $wb.Worksheets.Item("CUL")
# Row 4 (Leve Item ID 4650)
$ws.Range("H4").Value = 34655470
$ws.Range("I4").Value = 38654068
$ws.Range("K4").Value = 115962204
$ws.Range("M4").Value = -115962092

# Row 69 (Leve Item ID 12850)
$ws.Range("N69").ClearContents()
$ws.Range("H69").Value = 2500
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0

# Row 70 (Leve Item ID 12867)
$ws.Range("H70").Value = 324731.5
$ws.Range("J70").Value = 324731.5
$ws.Range("L70").Value = 974194.5
$ws.Range("N70").Value = -974824.5

# Row 72 (Leve Item ID 12850)
$ws.Range("N72").ClearContents()
$ws.Range("H72").Value = 2500
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0

# Row 73 (Leve Item ID 12867)
$ws.Range("H73").Value = 324731.5
$ws.Range("J73").Value = 324731.5
$ws.Range("L73").Value = 974194.5
$ws.Range("N73").Value = -976378.5

# Row 109 (Leve Item ID 27854)
$ws.Range("H109").Value = 380.8889
$ws.Range("I109").Value = 380.8889
$ws.Range("K109").Value = 1142.6667
$ws.Range("M109").Value = -102.6667

$ws = $wb.Worksheets.Item("GSM")
# Row 105 (Leve Item ID 18671)
$ws.Range("H105").Value = 150000
$ws.Range("J105").Value = 150000
$ws.Range("L105").Value = 150000
$ws.Range("N105").Value = -156988

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 5397.231
$ws.Range("I7").Value = 4980.7646
$ws.Range("K7").Value = 4980.7646
$ws.Range("M7").Value = -4868.7646

# Row 69 (Leve Item ID 10671)
$ws.Range("H69").Value = 58998
$ws.Range("J69").Value = 58998
$ws.Range("L69").Value = 58998
$ws.Range("N69").Value = -60620

# Row 72 (Leve Item ID 10671)
$ws.Range("H72").Value = 58998
$ws.Range("J72").Value = 58998
$ws.Range("L72").Value = 176994
$ws.Range("N72").Value = -185106

# Row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 8562.556
$ws.Range("J93").Value = 8717.75
$ws.Range("L93").Value = 8717.75
$ws.Range("N93").Value = -11213.75

# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 5397.231
$ws.Range("I126").Value = 4980.7646
$ws.Range("K126").Value = 14942.2938
$ws.Range("M126").Value = -12472.2938

# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 5430.6665
$ws.Range("I132").Value = 2836.6
$ws.Range("J132").Value = 10618.8
$ws.Range("K132").Value = 8509.799999999999
$ws.Range("L132").Value = 31856.4
$ws.Range("M132").Value = -5979.799999999999
$ws.Range("N132").Value = -36916.39999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 14 (Leve Item ID 2658)
$ws.Range("L14").ClearContents()
$ws.Range("H14").Value = 29998
$ws.Range("I14").Value = 29998
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 29998
$ws.Range("M14").Value = -29830
$ws.Range("N14").Value = 0

# Row 17 (Leve Item ID 3539)
$ws.Range("H17").Value = 2125
$ws.Range("I17").Value = 2125
$ws.Range("K17").Value = 2125
$ws.Range("M17").Value = -1953

# Row 38 (Leve Item ID 27990)
$ws.Range("H38").Value = 36528
$ws.Range("I38").Value = 23056
$ws.Range("K38").Value = 23056
$ws.Range("M38").Value = -22583

# Row 44 (Leve Item ID 2805)
$ws.Range("H44").Value = 41178.668
$ws.Range("J44").Value = 41178.668
$ws.Range("L44").Value = 41178.668
$ws.Range("N44").Value = -42286.668

# Row 54 (Leve Item ID 3413)
$ws.Range("L54").ClearContents()
$ws.Range("H54").Value = 26499.834
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 26499.834
$ws.Range("K54").Value = 0
$ws.Range("M54").Value = 26499.834
$ws.Range("N54").Value = -27539.834

# Row 81 (Leve Item ID 12596)
$ws.Range("H81").Value = 12349417
$ws.Range("I81").Value = 3106.3333
$ws.Range("J81").Value = 37042036
$ws.Range("K81").Value = 6212.6666
$ws.Range("L81").Value = 74084072
$ws.Range("M81").Value = -5151.6666
$ws.Range("N81").Value = -74086194

# Row 84 (Leve Item ID 12596)
$ws.Range("H84").Value = 12349417
$ws.Range("I84").Value = 3106.3333
$ws.Range("J84").Value = 37042036
$ws.Range("K84").Value = 31063.333
$ws.Range("L84").Value = 370420360
$ws.Range("M84").Value = -25759.333
$ws.Range("N84").Value = -370430968

# Row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 2185.4119
$ws.Range("J122").Value = 2537
$ws.Range("L122").Value = 7611
$ws.Range("N122").Value = -12511

# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 1415.6552
$ws.Range("I136").Value = 842.2
$ws.Range("K136").Value = 2526.6
$ws.Range("M136").Value = 23.39999999999964
